$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("GlobalConstantStringTable")
$ws3.Range("B1:B2").ColumnWidth = 13.375
